$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.193268
$ws.Range("H2").Value = 0.579804
$ws.Range("I2").Value = 0.1207017725010034
$ws.Range("J2").Value = 0.1207017725010034
$ws.Range("M2").Value = 19.827687
$ws.Range("N2").Value = 59.483061
$ws.Range("O2").Value = 0.1538389073329896
$ws.Range("P2").Value = 0.1538389073329896
$ws.Range("Q2").Value = 3.832057411116
$ws.Range("R2").Value = 34.488516700044
$ws.Range("S2").Value = 0.01856862879470946
$ws.Range("T2").Value = 0.01856862879470945
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.193268
$ws.Range("H3").Value = 0.579804
$ws.Range("I3").Value = 0.1207017725010034
$ws.Range("J3").Value = 0.1207017725010034
$ws.Range("O3").Value = 0.6604253914664442
$ws.Range("P3").Value = 0.6604253914664441
$ws.Range("Q3").Value = 16.450896978748
$ws.Range("R3").Value = 148.058072808732
$ws.Range("S3").Value = 0.07971451535466886
$ws.Range("T3").Value = 0.07971451535466884
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.193268
$ws.Range("H4").Value = 0.579804
$ws.Range("I4").Value = 0.1207017725010034
$ws.Range("J4").Value = 0.1207017725010034
$ws.Range("M4").Value = 23.93873833333333
$ws.Range("N4").Value = 71.816215
$ws.Range("O4").Value = 0.1857357012005663
$ws.Range("P4").Value = 0.1857357012005663
$ws.Range("Q4").Value = 4.626592080206667
$ws.Range("R4").Value = 41.63932872186
$ws.Range("S4").Value = 0.0224186283516251
$ws.Range("T4").Value = 0.02241862835162509
$ws.Range("I5").Value = 0.7607038292883184
$ws.Range("J5").Value = 0.7607038292883183
$ws.Range("M5").Value = 19.827687
$ws.Range("N5").Value = 59.483061
$ws.Range("O5").Value = 0.1538389073329896
$ws.Range("P5").Value = 0.1538389073329896
$ws.Range("Q5").Value = 24.150935701167
$ws.Range("R5").Value = 217.358421310503
$ws.Range("S5").Value = 0.117025845901736
$ws.Range("T5").Value = 0.1170258459017359
$ws.Range("I6").Value = 0.7607038292883184
$ws.Range("J6").Value = 0.7607038292883183
$ws.Range("O6").Value = 0.6604253914664442
$ws.Range("P6").Value = 0.6604253914664441
$ws.Range("S6").Value = 0.5023881242477608
$ws.Range("T6").Value = 0.5023881242477606
$ws.Range("I7").Value = 0.7607038292883184
$ws.Range("J7").Value = 0.7607038292883183
$ws.Range("M7").Value = 23.93873833333333
$ws.Range("N7").Value = 71.816215
$ws.Range("O7").Value = 0.1857357012005663
$ws.Range("P7").Value = 0.1857357012005663
$ws.Range("Q7").Value = 29.15836477827167
$ws.Range("R7").Value = 262.425283004445
$ws.Range("S7").Value = 0.1412898591388217
$ws.Range("T7").Value = 0.1412898591388217
$ws.Range("G8").Value = 0.1898936666666667
$ws.Range("H8").Value = 0.569681
$ws.Range("I8").Value = 0.1185943982106783
$ws.Range("J8").Value = 0.1185943982106783
$ws.Range("M8").Value = 19.827687
$ws.Range("N8").Value = 59.483061
$ws.Range("O8").Value = 0.1538389073329896
$ws.Range("P8").Value = 0.1538389073329896
$ws.Range("Q8").Value = 3.765152185949
$ws.Range("R8").Value = 33.886369673541
$ws.Range("S8").Value = 0.01824443263654421
$ws.Range("T8").Value = 0.01824443263654421
$ws.Range("G9").Value = 0.1898936666666667
$ws.Range("H9").Value = 0.569681
$ws.Range("I9").Value = 0.1185943982106783
$ws.Range("J9").Value = 0.1185943982106783
$ws.Range("O9").Value = 0.6604253914664442
$ws.Range("P9").Value = 0.6604253914664441
$ws.Range("Q9").Value = 16.16367503803033
$ws.Range("R9").Value = 145.473075342273
$ws.Range("S9").Value = 0.07832275186401458
$ws.Range("T9").Value = 0.07832275186401456
$ws.Range("G10").Value = 0.1898936666666667
$ws.Range("H10").Value = 0.569681
$ws.Range("I10").Value = 0.1185943982106783
$ws.Range("J10").Value = 0.1185943982106783
$ws.Range("M10").Value = 23.93873833333333
$ws.Range("N10").Value = 71.816215
$ws.Range("O10").Value = 0.1857357012005663
$ws.Range("P10").Value = 0.1857357012005663
$ws.Range("Q10").Value = 4.545814797490555
$ws.Range("R10").Value = 40.912333177415
$ws.Range("S10").Value = 0.02202721371011952
$ws.Range("T10").Value = 0.02202721371011952
